# Update "想去人数" (want-to-go count) figures in column F across the
# 展览 / 演出 / 本地生活 sheets and the combined 全部类型 sheet, matching the
# refreshed data snapshot from the gh-pages generation run.
$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5, 6).Value = 771
$ws.Cells.Item(6, 6).Value = 465
$ws.Cells.Item(8, 6).Value = 204
$ws.Cells.Item(11, 6).Value = 7195
$ws.Cells.Item(12, 6).Value = 78
$ws.Cells.Item(14, 6).Value = 1682
$ws.Cells.Item(20, 6).Value = 22
$ws.Cells.Item(21, 6).Value = 733
$ws.Cells.Item(22, 6).Value = 12
$ws.Cells.Item(24, 6).Value = 136
$ws.Cells.Item(25, 6).Value = 6
$ws.Cells.Item(26, 6).Value = 207
$ws.Cells.Item(30, 6).Value = 1061
$ws.Cells.Item(32, 6).Value = 85
$ws.Cells.Item(33, 6).Value = 2059
$ws.Cells.Item(34, 6).Value = 582
$ws.Cells.Item(35, 6).Value = 6
$ws.Cells.Item(36, 6).Value = 17
$ws.Cells.Item(38, 6).Value = 559

# 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value = 61
$ws.Cells.Item(6, 6).Value = 307

# 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 366

# 全部类型 (All types, combined listing)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 366
$ws.Cells.Item(6, 6).Value = 771
$ws.Cells.Item(8, 6).Value = 465
$ws.Cells.Item(10, 6).Value = 204
$ws.Cells.Item(13, 6).Value = 7195
$ws.Cells.Item(14, 6).Value = 78
$ws.Cells.Item(17, 6).Value = 1682
$ws.Cells.Item(22, 6).Value = 61
$ws.Cells.Item(24, 6).Value = 22
$ws.Cells.Item(26, 6).Value = 307
$ws.Cells.Item(28, 6).Value = 733
$ws.Cells.Item(29, 6).Value = 12
$ws.Cells.Item(31, 6).Value = 136
$ws.Cells.Item(33, 6).Value = 6
$ws.Cells.Item(36, 6).Value = 207
$ws.Cells.Item(40, 6).Value = 1061
$ws.Cells.Item(42, 6).Value = 85
$ws.Cells.Item(43, 6).Value = 2059
$ws.Cells.Item(44, 6).Value = 582
$ws.Cells.Item(45, 6).Value = 6
$ws.Cells.Item(46, 6).Value = 17
$ws.Cells.Item(48, 6).Value = 559
